$d = $word.ActiveDocument

# --- Edit 1: split the "define" word out into its own run with gramStart/gramEnd proofErr markers ---
$p4 = $d.Paragraphs.Item(4)
$rng1 = $d.Range($p4.Range.Start, $p4.Range.End)

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="78F4133F" w14:textId="77777777" w:rsidR="005A5940" w:rsidRPr="005A5940" w:rsidRDefault="005A5940" w:rsidP="005A5940"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="single" w:sz="2" w:space="0" w:color="E3E3E3"/><w:left w:val="single" w:sz="2" w:space="5" w:color="E3E3E3"/><w:bottom w:val="single" w:sz="2" w:space="0" w:color="E3E3E3"/><w:right w:val="single" w:sz="2" w:space="0" w:color="E3E3E3"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">Is file </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>mein</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">, main page ka structure </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>define</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>karein</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">. Yeh file header, footer, aur menu-sidebar components ko include </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>karega</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005A5940"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng1.InsertXML($xml1)

# --- Edit 2: append new paragraphs (Second commit ... tree listing) after the last paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng2 = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="2C48722B" w14:textId="47ED6712" w:rsidR="005A5940" w:rsidRDefault="005A5940" w:rsidP="005A5940"><w:r><w:t xml:space="preserve">      - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MenuSide.jsx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Second commit </w:t></w:r></w:p><w:p><w:r><w:t>pages</w:t></w:r></w:p><w:p><w:r><w:t>│   dashboard</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">│   │   </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dashboard.jsx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>│   users</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">│   │   </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userList.jsx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>│   permissions</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">│   │   </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>permissionList.jsx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>│   roles</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">│   │   </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>roleList.jsx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng2.InsertXML($xml2)

Write-Host "Done"
